$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format so numeric-looking strings
# ("226.86", "0.998", ...) are written back as text, matching the
# original inline-string cell type instead of being coerced to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '34.595.41'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '1.842.67'
$ws.Range('E3').Value = '  +3.98%  '
$ws.Range('E4').Value = '  -0.52%  '
$ws.Range('D5').Value = '226.86'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('D8').Value = '32.70'
$ws.Range('E8').Value = '  +5.48%  '
$ws.Range('E9').Value = '  +5.48%  '
$ws.Range('D10').Value = '0.0717'
$ws.Range('E10').Value = '  +9.91%  '
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '2.110.44'
$ws.Range('E12').Value = '  +4.05%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.843.16'
$ws.Range('E13').Value = '  +3.61%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '11.21'
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('E15').Value = '  +5.18%  '
$ws.Range('D16').Value = '34.616.96'
$ws.Range('E16').Value = '  +2.07%  '
$ws.Range('E17').Value = '  +4.22%  '
$ws.Range('D18').Value = '69.85'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('D19').Value = '254.09'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').Value = '0.0₃0809'
$ws.Range('E20').Value = '  +10.35%  '
$ws.Range('D21').Value = '11.35'
$ws.Range('E21').Value = '  +10.60%  '
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('D23').Value = '4.33'
$ws.Range('E23').Value = '  +3.91%  '
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').Value = '161.79'
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('D26').Value = '16.90'
$ws.Range('E26').Value = '  +3.77%  '
$ws.Range('D27').Value = '7.26'
$ws.Range('E27').Value = '  +4.75%  '
$ws.Range('D28').Value = '0.116'
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').Value = '0.0533'
$ws.Range('E30').Value = '  +4.51%  '
$ws.Range('E31').Value = '  +2.71%  '
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('D33').Value = '515.69'
$ws.Range('E33').Value = '  +886.14%  '
$ws.Range('E34').Value = '  +2.76%  '
$ws.Range('D35').Value = '1.95'
$ws.Range('E35').Value = '  +7.59%  '
$ws.Range('D36').Value = '1.445.93'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('E39').Value = '  +4.65%  '
$ws.Range('D40').Value = '0.982'
$ws.Range('E40').Value = '  +11.47%  '
$ws.Range('D41').Value = '83.23'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').Value = '2.81'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').Value = '2.37'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('D44').Value = '2.17'
$ws.Range('E44').Value = '  +6.42%  '
$ws.Range('D45').Value = '6.12'
$ws.Range('E45').Value = '  +6.34%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '12.61'
$ws.Range('E46').Value = '  +7.55%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '2.001.30'
$ws.Range('E47').Value = '  +3.82%  '
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').Value = '0.0493'
$ws.Range('E49').Value = '  -3.16%  '
$ws.Range('D50').Value = '106.55'
$ws.Range('E50').Value = '  +10.02%  '
$ws.Range('E51').Value = '  -0.18%  '

# Restore the default style so no stray formatting is left behind.
$ws.Range('D2:D51').Style = 'Normal'
